# "Added a few logos"
#
# Adds a new "logo" column (I) to the places sheet and fills it with a
# slug that identifies which logo image goes with each row. Most rows
# don't have a specific logo yet and fall back to the "168-sushi"
# placeholder; a handful of restaurants that were just given real
# logo artwork get their own slug (each of those restaurants spans two
# rows in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("I1").Value = "logo"

# 168 Sushi Buffet (rows 2-3)
$ws.Range("I2:I3").Value = "168-sushi"

# Bang Bang (row 4)
$ws.Range("I4").Value = "bang-bang"

# Banh Mi Boys (rows 5-6)
$ws.Range("I5:I6").Value = "banh-mi-boys"

# BeaverTails (rows 7-8)
$ws.Range("I7:I8").Value = "beavertails"

# Boba Boy (rows 9-10)
$ws.Range("I9:I10").Value = "boba-boy"

# Blaze Pizza (rows 11-12)
$ws.Range("I11:I12").Value = "blaze-pizza"

# Everyone else still uses the placeholder logo for now.
$ws.Range("I13:I63").Value = "168-sushi"

# Size the new column to fit its contents, like Excel does automatically
# after typing values into a fresh column.
$ws.Columns.Item(9).AutoFit()

# Leave the selection where it was when the edit was last saved.
$ws.Range("I13").Select()
